$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '51.096.65'
$ws.Range('E2').Value = '  -1.20%  '
$ws.Range('D3').Value = '2.911.53'
$ws.Range('E3').Value = '  -0.45%  '
$ws.Range('E4').Value = '  -0.13%  '
$ws.Range('D5').Value = '''369.05'
$ws.Range('E5').Value = '  +5.00%  '
$ws.Range('D6').Value = '''103.44'
$ws.Range('E6').Value = '  -2.67%  '
$ws.Range('E7').Value = '  -2.91%  '
$ws.Range('E8').Value = '  -0.10%  '
$ws.Range('E9').Value = '  -3.40%  '
$ws.Range('D10').Value = '''36.65'
$ws.Range('E10').Value = '  -2.73%  '
$ws.Range('E11').Value = '  +1.33%  '
$ws.Range('E12').Value = '  -2.06%  '
$ws.Range('D13').Value = '''18.37'
$ws.Range('E13').Value = '  -3.18%  '
$ws.Range('D14').Value = '3.371.43'
$ws.Range('E14').Value = '  -0.53%  '
$ws.Range('E15').Value = '  -3.60%  '
$ws.Range('D16').Value = '2.912.56'
$ws.Range('E16').Value = '  -0.50%  '
$ws.Range('D17').Value = '''0.944'
$ws.Range('E17').Value = '  -2.08%  '
$ws.Range('D18').Value = '51.028.70'
$ws.Range('E18').Value = '  -1.26%  '
$ws.Range('E19').Value = '  -4.65%  '
$ws.Range('D20').Value = '''7.22'
$ws.Range('E20').Value = '  -1.54%  '
$ws.Range('D21').Value = '''12.80'
$ws.Range('E21').Value = '  -4.42%  '
$ws.Range('D22').Value = '0.0₃0946'
$ws.Range('E22').Value = '  -1.51%  '
$ws.Range('D23').Value = '''68.26'
$ws.Range('E23').Value = '  -0.82%  '
$ws.Range('D24').Value = '''259.62'
$ws.Range('E24').Value = '  -0.84%  '
$ws.Range('D25').Value = '''2.68'
$ws.Range('E25').Value = '  -0.97%  '
$ws.Range('E26').Value = '  +0.53%  '
$ws.Range('D28').Value = '''25.71'
$ws.Range('E28').Value = '  -2.94%  '
$ws.Range('D29').Value = '''6.96'
$ws.Range('E29').Value = '  -5.63%  '
$ws.Range('D30').Value = '''0.102'
$ws.Range('E30').Value = '  -0.41%  '
$ws.Range('B31').Value = 'RenderToken'
$ws.Range('C31').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D31').Value = '''6.14'
$ws.Range('E31').Value = '  +3.94%  '
$ws.Range('B32').Value = 'Cosmos'
$ws.Range('C32').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D32').Value = '''9.91'
$ws.Range('E32').Value = '  -3.12%  '
$ws.Range('E33').Value = '  -1.53%  '
$ws.Range('D34').Value = '''34.69'
$ws.Range('E34').Value = '  -2.52%  '
$ws.Range('E35').Value = '  -0.30%  '
$ws.Range('E36').Value = '  +0.42%  '
$ws.Range('D37').Value = '''0.0421'
$ws.Range('E37').Value = '  -1.25%  '
$ws.Range('B38').Value = 'Stacks'
$ws.Range('C38').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D38').Value = '''2.68'
$ws.Range('E38').Value = '  +0.95%  '
$ws.Range('B39').Value = 'LidoDAOToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D39').Value = '''3.04'
$ws.Range('E39').Value = '  -3.18%  '
$ws.Range('D40').Value = '''17.09'
$ws.Range('E40').Value = '  -3.43%  '
$ws.Range('D41').Value = '''1.84'
$ws.Range('E41').Value = '  -5.74%  '
$ws.Range('E42').Value = '  -2.54%  '
$ws.Range('D43').Value = '''22.24'
$ws.Range('E43').Value = '  -1.24%  '
$ws.Range('E44').Value = '  -0.13%  '
$ws.Range('E45').Value = '  -3.10%  '
$ws.Range('D46').Value = '2.020.80'
$ws.Range('E46').Value = '  -3.82%  '
$ws.Range('D47').Value = '''2.31'
$ws.Range('E47').Value = '  -5.31%  '
$ws.Range('D48').Value = '''3.17'
$ws.Range('E48').Value = '  -4.47%  '
$ws.Range('D49').Value = '3.198.97'
$ws.Range('E49').Value = '  -0.35%  '
$ws.Range('D50').Value = '''0.239'
$ws.Range('E50').Value = '  +0.63%  '
$ws.Range('D51').Value = '''0.0312'
$ws.Range('E51').Value = '  -7.70%  '
